$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarterly scroll: columns D..M shift one quarter left, new quarter data lands in column M ---
# (period header row 8, publish-date row 9, and every financial data row 12-58)

$row8 = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
for ($i = 0; $i -lt $row8.Count; $i++) { $ws.Cells.Item(8, 4 + $i).Value = $row8[$i] }

$row9 = @("1399-10-30", "1401-02-25 (12)", "1400-05-20 (2)", "1400-08-29 (2)", "1400-11-02", "1402-02-27 (12)", "1401-05-01", "1401-08-30 (2)", "1401-11-19 (2)", "1402-02-27 (3)")
for ($i = 0; $i -lt $row9.Count; $i++) { $ws.Cells.Item(9, 4 + $i).Value = $row9[$i] }

$row12 = @(539172, 819625, 210560, 1421181, 882761, 1224333, 3409514, 3292993, 1708760, 1251834)
for ($i = 0; $i -lt $row12.Count; $i++) { $ws.Cells.Item(12, 4 + $i).Value = $row12[$i] }

$row13 = @(1301499, 2084468, 2130113, 951020, 3031778, 3166748, 3166748, 3713154, 3758124, 7051746)
for ($i = 0; $i -lt $row13.Count; $i++) { $ws.Cells.Item(13, 4 + $i).Value = $row13[$i] }

$row14 = @(668421, 243705, 200558, 167720, 90811, 153222, 167820, 217659, 400577, 593130)
for ($i = 0; $i -lt $row14.Count; $i++) { $ws.Cells.Item(14, 4 + $i).Value = $row14[$i] }

$row15 = @(1880417, 1936525, 2048790, 2230119, 2664833, 3055890, 2930414, 4216028, 4348400, 5147002)
for ($i = 0; $i -lt $row15.Count; $i++) { $ws.Cells.Item(15, 4 + $i).Value = $row15[$i] }

$row16 = @(236120, 301162, 274095, 281036, 237817, 403210, 317726, 194236, 630328, 1039192)
for ($i = 0; $i -lt $row16.Count; $i++) { $ws.Cells.Item(16, 4 + $i).Value = $row16[$i] }

$row17 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row17.Count; $i++) { $ws.Cells.Item(17, 4 + $i).Value = $row17[$i] }

$row18 = @(4625629, 5385485, 4864116, 5051076, 6908000, 8003403, 9992222, 11634070, 10846189, 15082904)
for ($i = 0; $i -lt $row18.Count; $i++) { $ws.Cells.Item(18, 4 + $i).Value = $row18[$i] }

$row19 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row19.Count; $i++) { $ws.Cells.Item(19, 4 + $i).Value = $row19[$i] }

$row20 = @(6258, 6258, 678147, 552791, 628811, 566248, 601477, 566248, 922755, 1238176)
for ($i = 0; $i -lt $row20.Count; $i++) { $ws.Cells.Item(20, 4 + $i).Value = $row20[$i] }

$row21 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row21.Count; $i++) { $ws.Cells.Item(21, 4 + $i).Value = $row21[$i] }

$row22 = @(517135, 531044, 472533, 528566, 524645, 521954, 512973, 510163, 633696, 684469)
for ($i = 0; $i -lt $row22.Count; $i++) { $ws.Cells.Item(22, 4 + $i).Value = $row22[$i] }

$row23 = @(58441, 58441, 58441, 58441, 58441, 58441, 58441, 58441, 60391, 61867)
for ($i = 0; $i -lt $row23.Count; $i++) { $ws.Cells.Item(23, 4 + $i).Value = $row23[$i] }

$row24 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $row24.Count; $i++) { $ws.Cells.Item(24, 4 + $i).Value = $row24[$i] }

$row25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row25.Count; $i++) { $ws.Cells.Item(25, 4 + $i).Value = $row25[$i] }

$row26 = @(581834, 595743, 1209121, 1139798, 1211897, 1146643, 1172891, 1134852, 1616842, 1984512)
for ($i = 0; $i -lt $row26.Count; $i++) { $ws.Cells.Item(26, 4 + $i).Value = $row26[$i] }

$row27 = @(5207463, 5981228, 6073237, 6190874, 8119897, 9150046, 11165113, 12768922, 12463031, 17067416)
for ($i = 0; $i -lt $row27.Count; $i++) { $ws.Cells.Item(27, 4 + $i).Value = $row27[$i] }

$row29 = @(492741, 448646, 531435, 631991, 830963, 1083385, 913570, 2560323, 925703, 1891976)
for ($i = 0; $i -lt $row29.Count; $i++) { $ws.Cells.Item(29, 4 + $i).Value = $row29[$i] }

$row30 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $row30.Count; $i++) { $ws.Cells.Item(30, 4 + $i).Value = $row30[$i] }

$row31 = @(366025, 477357, 195113, 1007589, 936668, 681268, 1105663, 1072856, 714195, 1017875)
for ($i = 0; $i -lt $row31.Count; $i++) { $ws.Cells.Item(31, 4 + $i).Value = $row31[$i] }

$row32 = @(418593, 457817, 0, 740281, 734701, 758215, 645142, 1400908, 1655015, 1838730)
for ($i = 0; $i -lt $row32.Count; $i++) { $ws.Cells.Item(32, 4 + $i).Value = $row32[$i] }

$row33 = @(143441, 116424, 2180228, 150483, 126335, 125994, 5025972, 1828709, 179836, 147480)
for ($i = 0; $i -lt $row33.Count; $i++) { $ws.Cells.Item(33, 4 + $i).Value = $row33[$i] }

$row34 = @(34320, 34320, 34320, 36470, 27890, 34320, 34320, 19389, 17160, 0)
for ($i = 0; $i -lt $row34.Count; $i++) { $ws.Cells.Item(34, 4 + $i).Value = $row34[$i] }

$row35 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row35.Count; $i++) { $ws.Cells.Item(35, 4 + $i).Value = $row35[$i] }

$row36 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row36.Count; $i++) { $ws.Cells.Item(36, 4 + $i).Value = $row36[$i] }

$row37 = @(1455120, 1534564, 2941096, 2566814, 2656557, 2683182, 7724667, 6882185, 3491909, 4896061)
for ($i = 0; $i -lt $row37.Count; $i++) { $ws.Cells.Item(37, 4 + $i).Value = $row37[$i] }

$row38 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row38.Count; $i++) { $ws.Cells.Item(38, 4 + $i).Value = $row38[$i] }

$row39 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $row39.Count; $i++) { $ws.Cells.Item(39, 4 + $i).Value = $row39[$i] }

$row40 = @(51480, 34320, 34320, 17160, 17160, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row40.Count; $i++) { $ws.Cells.Item(40, 4 + $i).Value = $row40[$i] }

$row41 = @(241238, 242009, 325502, 318100, 336577, 353916, 532744, 546069, 554566, 614928)
for ($i = 0; $i -lt $row41.Count; $i++) { $ws.Cells.Item(41, 4 + $i).Value = $row41[$i] }

$row42 = @(292718, 276329, 359822, 335260, 353737, 353916, 532744, 546069, 554566, 614928)
for ($i = 0; $i -lt $row42.Count; $i++) { $ws.Cells.Item(42, 4 + $i).Value = $row42[$i] }

$row43 = @(1747838, 1810893, 3300918, 2902074, 3010294, 3037098, 8257411, 7428254, 4046475, 5510989)
for ($i = 0; $i -lt $row43.Count; $i++) { $ws.Cells.Item(43, 4 + $i).Value = $row43[$i] }

$row45 = @(650000, 650000, 650000, 650000, 650000, 650000, 650000, 650000, 650000, 650000)
for ($i = 0; $i -lt $row45.Count; $i++) { $ws.Cells.Item(45, 4 + $i).Value = $row45[$i] }

$row46 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row46.Count; $i++) { $ws.Cells.Item(46, 4 + $i).Value = $row46[$i] }

$row47 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row47.Count; $i++) { $ws.Cells.Item(47, 4 + $i).Value = $row47[$i] }

$row48 = @(0, -66171, 0, -13728, -13728, -58111, -58111, -13, -13, -3440)
for ($i = 0; $i -lt $row48.Count; $i++) { $ws.Cells.Item(48, 4 + $i).Value = $row48[$i] }

$row49 = @(0, 7168, 0, 0, 0, 555, 555, 8458, 8458, 4785)
for ($i = 0; $i -lt $row49.Count; $i++) { $ws.Cells.Item(49, 4 + $i).Value = $row49[$i] }

$row50 = @(65000, 65000, 65000, 65000, 65000, 65000, 65000, 65000, 65000, 65000)
for ($i = 0; $i -lt $row50.Count; $i++) { $ws.Cells.Item(50, 4 + $i).Value = $row50[$i] }

$row51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row51.Count; $i++) { $ws.Cells.Item(51, 4 + $i).Value = $row51[$i] }

$row52 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $row52.Count; $i++) { $ws.Cells.Item(52, 4 + $i).Value = $row52[$i] }

$row53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row53.Count; $i++) { $ws.Cells.Item(53, 4 + $i).Value = $row53[$i] }

$row54 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
for ($i = 0; $i -lt $row54.Count; $i++) { $ws.Cells.Item(54, 4 + $i).Value = $row54[$i] }

$row55 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $row55.Count; $i++) { $ws.Cells.Item(55, 4 + $i).Value = $row55[$i] }

$row56 = @(2744625, 3514338, 2057319, 2587528, 4408331, 5455504, 2250258, 4617223, 7693111, 10840082)
for ($i = 0; $i -lt $row56.Count; $i++) { $ws.Cells.Item(56, 4 + $i).Value = $row56[$i] }

$row57 = @(3459625, 4170335, 2772319, 3288800, 5109603, 6112948, 2907702, 5340668, 8416556, 11556427)
for ($i = 0; $i -lt $row57.Count; $i++) { $ws.Cells.Item(57, 4 + $i).Value = $row57[$i] }

$row58 = @(5207463, 5981228, 6073237, 6190874, 8119897, 9150046, 11165113, 12768922, 12463031, 17067416)
for ($i = 0; $i -lt $row58.Count; $i++) { $ws.Cells.Item(58, 4 + $i).Value = $row58[$i] }

# --- Column width adjustments (widths also scroll left by one column) ---
$ws.Range("E1").ColumnWidth = 30.166666666666668
$ws.Range("F1").ColumnWidth = 28.166666666666668
$ws.Range("I1").ColumnWidth = 30.166666666666668
$ws.Range("J1").ColumnWidth = 28.166666666666668
$ws.Range("M1").ColumnWidth = 30.166666666666668
